$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-coerced into numbers by Excel (losing significant trailing zeros),
# matching the original inlineStr text formatting of column D.
$textCells = @("D5", "D6", "D10", "D13", "D14", "D16", "D20", "D21", "D22", "D24", "D26", "D29", "D30", "D31", "D34", "D36", "D37", "D40", "D42", "D44", "D46", "D49", "D50")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply updated coin price / 1h volume values (and the TRON / ShibaInu row swap)
$ws.Range("D2").Value = '64.282.37'
$ws.Range("E2").Value = '  +0.81%  '
$ws.Range("D3").Value = '3.498.26'
$ws.Range("E3").Value = '  +0.48%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '585.82'
$ws.Range("E5").Value = '  +0.32%  '
$ws.Range("D6").Value = '134.22'
$ws.Range("E6").Value = '  +2.31%  '
$ws.Range("E8").Value = '  +0.76%  '
$ws.Range("E9").Value = '  +1.43%  '
$ws.Range("D10").Value = '7.27'
$ws.Range("E10").Value = '  +0.82%  '
$ws.Range("E11").Value = '  +1.69%  '
$ws.Range("D12").Value = '4.096.37'
$ws.Range("E12").Value = '  +0.58%  '
$ws.Range("B13").Value = 'ShibaInu'
$ws.Range("C13").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D13").Value = '0.0000182'
$ws.Range("E13").Value = '  +3.12%  '
$ws.Range("B14").Value = 'TRON'
$ws.Range("C14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D14").Value = '0.120'
$ws.Range("E14").Value = '  +1.18%  '
$ws.Range("D15").Value = '3.499.50'
$ws.Range("E15").Value = '  +0.98%  '
$ws.Range("D16").Value = '26.07'
$ws.Range("E16").Value = '  -4.81%  '
$ws.Range("D17").Value = '64.302.88'
$ws.Range("E17").Value = '  +0.80%  '
$ws.Range("E18").Value = '  +0.62%  '
$ws.Range("E19").Value = '  +1.67%  '
$ws.Range("D20").Value = '13.72'
$ws.Range("D21").Value = '393.55'
$ws.Range("E21").Value = '  +2.63%  '
$ws.Range("D22").Value = '0.571'
$ws.Range("E22").Value = '  -0.81%  '
$ws.Range("D23").Value = '3.640.24'
$ws.Range("E23").Value = '  +0.52%  '
$ws.Range("D24").Value = '74.21'
$ws.Range("E24").Value = '  +1.30%  '
$ws.Range("E25").Value = '  -0.16%  '
$ws.Range("D26").Value = '5.65'
$ws.Range("E26").Value = '  -1.10%  '
$ws.Range("E27").Value = '  +1.00%  '
$ws.Range("E28").Value = '  -0.07%  '
$ws.Range("D29").Value = '7.44'
$ws.Range("E29").Value = '  -1.58%  '
$ws.Range("D30").Value = '1.51'
$ws.Range("E30").Value = '  -4.72%  '
$ws.Range("D31").Value = '8.30'
$ws.Range("E31").Value = '  +0.17%  '
$ws.Range("E32").Value = '  -0.11%  '
$ws.Range("D33").Value = '3.517.97'
$ws.Range("E33").Value = '  +0.80%  '
$ws.Range("D34").Value = '0.151'
$ws.Range("E34").Value = '  +4.17%  '
$ws.Range("D36").Value = '23.47'
$ws.Range("E36").Value = '  +0.18%  '
$ws.Range("D37").Value = '5.19'
$ws.Range("E37").Value = '  -2.45%  '
$ws.Range("E38").Value = '  +0.29%  '
$ws.Range("E39").Value = '  -0.84%  '
$ws.Range("D40").Value = '161.43'
$ws.Range("E40").Value = '  +1.14%  '
$ws.Range("E41").Value = '  -1.59%  '
$ws.Range("D42").Value = '0.807'
$ws.Range("E42").Value = '  -0.17%  '
$ws.Range("E43").Value = '  +0.05%  '
$ws.Range("D44").Value = '25.10'
$ws.Range("E44").Value = '  -4.72%  '
$ws.Range("E45").Value = '  +0.29%  '
$ws.Range("D46").Value = '1.17'
$ws.Range("E46").Value = '  -3.09%  '
$ws.Range("E47").Value = '  +1.95%  '
$ws.Range("D48").Value = '2.466.92'
$ws.Range("E48").Value = '  +2.05%  '
$ws.Range("D49").Value = '6.78'
$ws.Range("E49").Value = '  -0.86%  '
$ws.Range("D50").Value = '0.896'
$ws.Range("E50").Value = '  -0.78%  '
$ws.Range("E51").Value = '  -1.23%  '
